# Applies the commit "Add files via upload" changes to red_ahash.xlsx sheet:
#  - resets several "processing time" measurements in column C to new timings
#  - appends 20 new rows (red_fred_1.jpg .. red_fred_20.jpg) with their hash data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column C ("Время обработки") values for existing rows ---
$ws.Range("C2").Value = 0
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0.015616
$ws.Range("C13").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("C21").Value = 0
$ws.Range("C26").Value = 0.015624
$ws.Range("C27").Value = 0
$ws.Range("C33").Value = 0
$ws.Range("C38").Value = 0
$ws.Range("C40").Value = 0.015624
$ws.Range("C44").Value = 0
$ws.Range("C48").Value = 0
$ws.Range("C49").Value = 0
$ws.Range("C50").Value = 0.015627
$ws.Range("C51").Value = 0
$ws.Range("C53").Value = 0
$ws.Range("C55").Value = 0
$ws.Range("C56").Value = 0.015626
$ws.Range("C57").Value = 0
$ws.Range("C58").Value = 0
$ws.Range("C60").Value = 0
$ws.Range("C62").Value = 0.015624
$ws.Range("C63").Value = 0
$ws.Range("C64").Value = 0

# --- Append 20 new rows (65-84) for red_fred_1.jpg .. red_fred_20.jpg ---
$hashRange = $ws.Range("B65:B84")
$hashRange.NumberFormat = "@"  # force text so the numeric-looking hash string is not coerced to a number

$ws.Range("A65").Value = "red_fred_1.jpg"
$ws.Range("B65").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 0

$ws.Range("A66").Value = "red_fred_2.jpg"
$ws.Range("B66").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C66").Value = 0.015624
$ws.Range("D66").Value = 0

$ws.Range("A67").Value = "red_fred_3.jpg"
$ws.Range("B67").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C67").Value = 0
$ws.Range("D67").Value = 0

$ws.Range("A68").Value = "red_fred_4.jpg"
$ws.Range("B68").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C68").Value = 0
$ws.Range("D68").Value = 0

$ws.Range("A69").Value = "red_fred_5.jpg"
$ws.Range("B69").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C69").Value = 0
$ws.Range("D69").Value = 0

$ws.Range("A70").Value = "red_fred_6.jpg"
$ws.Range("B70").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 0

$ws.Range("A71").Value = "red_fred_7.jpg"
$ws.Range("B71").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 0

$ws.Range("A72").Value = "red_fred_8.jpg"
$ws.Range("B72").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C72").Value = 0.015624
$ws.Range("D72").Value = 0

$ws.Range("A73").Value = "red_fred_9.jpg"
$ws.Range("B73").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C73").Value = 0
$ws.Range("D73").Value = 0

$ws.Range("A74").Value = "red_fred_10.jpg"
$ws.Range("B74").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C74").Value = 0
$ws.Range("D74").Value = 0

$ws.Range("A75").Value = "red_fred_11.jpg"
$ws.Range("B75").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C75").Value = 0
$ws.Range("D75").Value = 0

$ws.Range("A76").Value = "red_fred_12.jpg"
$ws.Range("B76").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C76").Value = 0
$ws.Range("D76").Value = 0

$ws.Range("A77").Value = "red_fred_13.jpg"
$ws.Range("B77").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C77").Value = 0.015625
$ws.Range("D77").Value = 0

$ws.Range("A78").Value = "red_fred_14.jpg"
$ws.Range("B78").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C78").Value = 0
$ws.Range("D78").Value = 0

$ws.Range("A79").Value = "red_fred_15.jpg"
$ws.Range("B79").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 0

$ws.Range("A80").Value = "red_fred_16.jpg"
$ws.Range("B80").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 0

$ws.Range("A81").Value = "red_fred_17.jpg"
$ws.Range("B81").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 0

$ws.Range("A82").Value = "red_fred_18.jpg"
$ws.Range("B82").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C82").Value = 0
$ws.Range("D82").Value = 0

$ws.Range("A83").Value = "red_fred_19.jpg"
$ws.Range("B83").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C83").Value = 0.015625
$ws.Range("D83").Value = 0

$ws.Range("A84").Value = "red_fred_20.jpg"
$ws.Range("B84").Value = "1110001111000001100000011000000110000001100000011100001111100011"
$ws.Range("C84").Value = 0
$ws.Range("D84").Value = 0

$hashRange.ClearFormats()  # drop the temporary text format so cells keep default styling
